$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column A spacer, narrow width (closest achievable to 4.140625 given engine's
# pixel-grid rounding on ColumnWidth -> OOXML width conversion)
$ws.Columns.Item(1).ColumnWidth = 3.25

# Fill in the remaining Stat Breakdown data (F:K) for rows 18-37
$ws.Range("F18").Value = 'Critical Strike'
$ws.Range("I18").Value = 'Fortified Leech'
$ws.Range("J18").Value = 'Celestial Guidance'
$ws.Range("F19").Value = 'Versatility'
$ws.Range("G19").Value = 'Stone'
$ws.Range("I19").Value = 'Fortified Avoidance'
$ws.Range("J19").Value = 'Celestial Guidance'
$ws.Range("F20").Value = 'Versatility'
$ws.Range("G20").Value = 'Shadowcore Oil'
$ws.Range("I20").Value = 'Fortified Leech'
$ws.Range("J20").Value = 'Celestial Guidance'
$ws.Range("F21").Value = 'Mastery'
$ws.Range("G21").Value = 'Shadowcore Oil'
$ws.Range("I21").Value = 'Fortified Leech'
$ws.Range("J21").Value = 'Celestial Guidance'
$ws.Range("F22").Value = 'Haste'
$ws.Range("G22").Value = 'Oil & Shapening Stone'
$ws.Range("I22").Value = 'Fortified Speed'
$ws.Range("J22").Value = 'Silful Revelation'
$ws.Range("I23").Value = 'Fortified Leech'
$ws.Range("J23").Value = 'Celestial Guidance'
$ws.Range("F24").Value = 'Haste'
$ws.Range("I24").Value = 'Fortified Leech'
$ws.Range("J24").Value = 'Celestial Guidance'
$ws.Range("F25").Value = 'Haste'
$ws.Range("I25").Value = 'Fortified Speed'
$ws.Range("J25").Value = 'Celestial Guidance'
$ws.Range("F26").Value = 'Critical Strike'
$ws.Range("G26").Value = 'Stone'
$ws.Range("I26").Value = 'Fortified Leech'
$ws.Range("J26").Value = 'Silful Revelation'
$ws.Range("K26").Value = 'Celestial Guidance'
$ws.Range("F27").Value = 'Versatility'
$ws.Range("G27").Value = 'Stone'
$ws.Range("I27").Value = 'Fortified Avoidance'
$ws.Range("J27").Value = 'Silful Revelation'
$ws.Range("K27").Value = 'Celestial Guidance'
$ws.Range("F28").Value = 'Versatility'
$ws.Range("G28").Value = 'Stone'
$ws.Range("I28").Value = 'Fortified Avoidance'
$ws.Range("J28").Value = 'Silful Revelation'
$ws.Range("K28").Value = 'Celestial Guidance'
$ws.Range("F29").Value = 'Critical Strike'
$ws.Range("I29").Value = 'Fortified Leech'
$ws.Range("J29").Value = 'Celestial Guidance'
$ws.Range("F30").Value = 'Haste'
$ws.Range("G30").Value = 'N/A'
$ws.Range("I30").Value = 'Fortified Avoidance'
$ws.Range("J30").Value = 'Silful Revelation'
$ws.Range("K30").Value = 'Celestial Guidance'
$ws.Range("I31").Value = 'Fortified Leech'
$ws.Range("J31").Value = 'Celestial Guidance'
$ws.Range("F32").Value = 'Haste'
$ws.Range("I32").Value = 'Fortified Speed'
$ws.Range("J32").Value = 'Celestial Guidance'
$ws.Range("F33").Value = 'Haste'
$ws.Range("I33").Value = '30 Stamina'
$ws.Range("J33").Value = 'Celestial Guidance'
$ws.Range("F34").Value = 'Haste'
$ws.Range("I34").Value = 'Fortified Speed'
$ws.Range("J34").Value = 'Celestial Guidance'
$ws.Range("F35").Value = 'Critical Strike'
$ws.Range("G35").Value = 'Stone'
$ws.Range("I35").Value = 'Fortified Avoidance'
$ws.Range("J35").Value = 'Silful Revelation'
$ws.Range("F36").Value = 'Mastery'
$ws.Range("G36").Value = 'Stone'
$ws.Range("I36").Value = 'Fortified Avoidance'
$ws.Range("J36").Value = 'Silful Revelation'
$ws.Range("K36").Value = 'Celestial Guidance'
$ws.Range("F37").Value = 'Versatility'
$ws.Range("G37").Value = 'Stone'
$ws.Range("I37").Value = 'Fortified Leech'
$ws.Range("J37").Value = 'Celestial Guidance'

# Restore the active selection to match the final edit position
$ws.Range("G37").Select()
